$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: to write a cell as a genuine *text* shared-string value
# (even when the text looks like a number, e.g. "2", "1", "0") without
# disturbing the cell's existing style, we stage the text in a scratch
# cell via a text formula (="...") , copy it, and Paste-Special "values
# only" into the destination. This preserves the destination's existing
# number format / style index (no new cellXfs entries get created) while
# guaranteeing the stored type is a string.
$scratch = $ws.Range("ZZ1")

function Set-TextCell($rangeAddr, $text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
}

# ---- Row 1 (header) ----
# A1 already holds "2" (unchanged by the edit) - leave it untouched so it
# keeps occupying the first shared-string slot.
Set-TextCell "B1" "1"
Set-TextCell "C1" "count"
Set-TextCell "D1" "0"

# ---- Row 2 (B-AP) ----
Set-TextCell "A2" "B-AP"
Set-TextCell "B2" "A, Np"
$ws.Range("C2").Value = 2871
Set-TextCell "D2" "nhiều, hơn, khác, cùng, nhất, nhỏ, lớn, gần, đầy, vui"

# ---- Row 3 (B-NP) ----
Set-TextCell "A3" "B-NP"
Set-TextCell "B3" "N, P, Nc, M, Np, L, Nu, Ny, FW, A"
$ws.Range("C3").Value = 17423
Set-TextCell "D3" "một, những, người, mình, tôi, đó, khi, chị, các, năm"

# ---- Row 4 (B-PP) ----
Set-TextCell "A4" "B-PP"
Set-TextCell "B4" "E"
$ws.Range("C4").Value = 2818
Set-TextCell "D4" "của, trong, với, cho, ở, để, về, từ, đến, trên"

# ---- Row 5 (B-VP) ----
Set-TextCell "A5" "B-VP"
Set-TextCell "B5" "V"
$ws.Range("C5").Value = 8715
Set-TextCell "D5" "có, là, đi, được, ra, lại, làm, nói, phải, biết"

# ---- Row 6 (I-NP) ----
Set-TextCell "A6" "I-NP"
Set-TextCell "B6" "Np, N, M, Ny, CH, V, A, FW, Nc, E"
$ws.Range("C6").Value = 827
Set-TextCell "D6" "HCM, Văn, Thị, Trâm, Nam, Trường Sơn, Lao Bảo, VN, Bắc, Đằng"

# ---- Row 7 (O) ----
Set-TextCell "A7" "O"
Set-TextCell "B7" "CH, R, C, T, X, I, Z, M"
$ws.Range("C7").Value = 11052
Set-TextCell "D7" ',, ., ", đã, và, không, ..., cũng, được, :'

# Clean up the scratch cell so it doesn't widen the sheet's used range.
$scratch.Clear()
